$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "弱电" (2nd sheet): a new row was inserted at row 50, shifting the
# existing rows 50-68 down to 51-69. The new row 50 duplicates the A/B/D/E/F/G
# values of the (then) following hydrant row but leaves column C blank.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Rows.Item(50).Insert()
$ws2.Rows.Item(50).RowHeight = 27

$ws2.Range("A50").Value = "室内消火栓平面"
$ws2.Range("B50").Value = "块:室内消火栓平面"
$ws2.Range("D50").Value = "E-BFAS610"
$ws2.Range("E50").Value = "块:E-BFAS610"
$ws2.Range("F50").Value = 0
$ws2.Range("G50").Value = "E-FAS-DEVC"

# ---------------------------------------------------------------------------
# View state: the workbook was left with the "弱电" sheet active, scrolled so
# that row 43 is near the top, with C50 selected.
# ---------------------------------------------------------------------------
$ws2.Activate()
$excel.ActiveWindow.ScrollRow = 43
$excel.ActiveWindow.ScrollColumn = 1
$ws2.Range("C50").Select()
